$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ("Movie Ratings") updates ---

# Analysis Completed? : "No" -> "Yes"
$ws.Range("J3").Value2 = "Yes"

# Date Added (H3) gets the same explicit date format used by H2, consolidating
# the duplicate "yyyy/mm/dd" number format onto the canonical one.
$ws.Range("H3").NumberFormat = "yyyy/mm/dd"

# Logo (K3) hyperlink, same styling convention as K2/L2/I3 (Arial 10, blue font, no underline).
$ws.Range("K3").Value2 = "https://github.com/kjeshang/KunalMavenAnalyticsDataPlayground/blob/main/Movie_Ratings/Logo.png?raw=true"
$ws.Hyperlinks.Add($ws.Range("K3"), "https://github.com/kjeshang/KunalMavenAnalyticsDataPlayground/blob/main/Movie_Ratings/Logo.png?raw=true", [Type]::Missing, [Type]::Missing, "https://github.com/kjeshang/KunalMavenAnalyticsDataPlayground/blob/main/Movie_Ratings/Logo.png?raw=true")
$ws.Range("K3").Font.Name = "Arial"
$ws.Range("K3").Font.Size = 10
$ws.Range("K3").Font.Color = 16711680
$ws.Range("K3").Font.Underline = 0

# Project Link (L3) hyperlink.
$ws.Range("L3").Value2 = "https://github.com/kjeshang/KunalMavenAnalyticsDataPlayground/tree/main/Movie_Ratings"
$ws.Hyperlinks.Add($ws.Range("L3"), "https://github.com/kjeshang/KunalMavenAnalyticsDataPlayground/tree/main/Movie_Ratings", [Type]::Missing, [Type]::Missing, "https://github.com/kjeshang/KunalMavenAnalyticsDataPlayground/tree/main/Movie_Ratings")
$ws.Range("L3").Font.Name = "Arial"
$ws.Range("L3").Font.Size = 10
$ws.Range("L3").Font.Color = 16711680
$ws.Range("L3").Font.Underline = 0

# Move the selection cursor below the table, like the source edit did.
$ws.Range("A4").Select() | Out-Null
